# Weekly update: a new price record for the current week is inserted at the
# top of the Apio (Feria Lagunitas de Puerto Montt) date-ordered data block,
# pushing all the existing records (old rows 258-353) down by one row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 258; everything from 258 downward shifts to 259+
$ws.Rows(258).Insert()

# Populate the new row with this week's record
$ws.Range("A258").Value = 4
$ws.Range("B258").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C258").Value = "Los Lagos"
$ws.Range("D258").Value = 44900
$ws.Range("E258").Value = 10
$ws.Range("F258").Value = 100112017
$ws.Range("G258").Value = "Apio"
$ws.Range("H258").Value = "Americana (o)"
$ws.Range("I258").Value = "Primera"
$ws.Range("J258").Value = 25
$ws.Range("K258").Value = 14000
$ws.Range("L258").Value = 14000
$ws.Range("M258").Value = 14000
$ws.Range("N258").Value = "`$/docena de matas"
$ws.Range("O258").Value = "Región de Coquimbo"
$ws.Range("P258").Value = 2333
$ws.Range("Q258").Value = 6
$ws.Range("R258").Value = "Hortaliza"
